$wb = $excel.ActiveWorkbook

# --- Rename the third sheet (Лист3 -> Radio) and populate it with the
#     "Radio" data-transmission timing calculations ---------------------
$ws = $wb.Worksheets.Item(3)
$ws.Name = "Radio"

# Column A width (chars) - closest reachable value to the recorded
# 14.7109375 stored width given this engine's column-width quantisation.
$ws.Columns.Item(1).ColumnWidth = 13.86

# Row 3 - Bitrate
$ws.Range("A3").Value = "Bitrate"
$ws.Range("C3").Value = "Bit/s"
$ws.Range("B3").Value = 10000

# Row 4 - Bit duration
$ws.Range("A4").Value = "Bit duration"
$ws.Range("C4").Value = "uS"
$ws.Range("B4").Formula = "=1000000/B3"

# Row 5 - Byte duration
$ws.Range("A5").Value = "Byte duration"
$ws.Range("C5").Value = "uS"
$ws.Range("B5").Formula = "=B4*8"

# Row 8 - "Experimental" section header
$ws.Range("A8").Value = "Experimental"
$ws.Range("A8").Style = "Accent1"

# Row 9 - Packet duration (note: C9 entered before A9 originally)
$ws.Range("C9").Value = "mS"
$ws.Range("A9").Value = "Packet duration"
$ws.Range("B9").Value = 24
$ws.Range("A9").Style = "Good"
$ws.Range("B9").Style = "Good"

# Row 10 - Receive duration (= packet duration * 2)
$ws.Range("A10").Value = "Receive duration"
$ws.Range("C10").Value = "ms"
$ws.Range("D10").Value = "*2"
$ws.Range("B10").Formula = "=B9*2"
$ws.Range("A10").Style = "Calculation"
$ws.Range("B10").Style = "Calculation"

# Row 12 - RX_OFF duration
$ws.Range("A12").Value = "RX_OFF duration"
$ws.Range("C12").Value = "ms"
$ws.Range("B12").Value = 216
$ws.Range("A12").Style = "Good"
$ws.Range("B12").Style = "Good"

# Row 13 - RX_ON duration
$ws.Range("A13").Value = "RX_ON duration"
$ws.Range("C13").Value = "ms"
$ws.Range("B13").Value = 54
$ws.Range("A13").Style = "Good"
$ws.Range("B13").Style = "Good"

# Row 14 - ON/OFF ratio
$ws.Range("A14").Value = "ON/OFF ratio"
$ws.Range("C14").Value = "%"
$ws.Range("B14").Formula = "=100*B13/(B13+B12)"
$ws.Range("A14").Style = "Calculation"
$ws.Range("B14").Style = "Calculation"
$ws.Range("B14").NumberFormat = "0.0"

# Selection + sheet activation (also clears tabSelected on the previously
# active sheet and sets bookViews/activeTab to this sheet).
$ws.Range("B12").Select()
$ws.Activate()
